$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header U1 from sd_80plus to sd_80plus_hosp
$ws.Range("U1").Value = "sd_80plus_hosp"

# Add tested value for row 30 (Apr 3)
$ws.Range("B30").Value = 1025

# Row 31 (Apr 4)
$ws.Range("A31").Value = 43925
$ws.Range("B31").Value = 807
$ws.Range("C31").Value = 1209
$ws.Range("D31").Value = 10
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 14
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = 209
$ws.Range("I31").Value = 13
$ws.Range("J31").Value = 251
$ws.Range("K31").Value = 26
$ws.Range("L31").Value = 212
$ws.Range("M31").Value = 30
$ws.Range("N31").Value = 202
$ws.Range("O31").Value = 41
$ws.Range("P31").Value = 152
$ws.Range("Q31").Value = 43
$ws.Range("R31").Value = 91
$ws.Range("S31").Value = 36
$ws.Range("T31").Value = 64
$ws.Range("U31").Value = 36
$ws.Range("V31").Value = 4
$ws.Range("W31").Value = 1
$ws.Range("X31").Value = 570
$ws.Range("Y31").Value = 630
$ws.Range("Z31").Value = 9
$ws.Range("AA31").Value = 228
$ws.Range("AB31").Value = 89
$ws.Range("AC31").Value = 18
$ws.Range("AD31").Value = 39
$ws.Range("AE31").Value = 86
$ws.Range("AF31").Value = 4
$ws.Range("AG31").Value = 7
$ws.Range("AH31").Value = 67
$ws.Range("AI31").Value = 27
$ws.Range("AJ31").Value = 27
$ws.Range("AK31").Value = 3
$ws.Range("AL31").Value = 18
$ws.Range("AM31").Value = 13
$ws.Range("AN31").Value = 18
$ws.Range("AO31").Value = 29
$ws.Range("AP31").Value = 13
$ws.Range("AQ31").Value = 649
$ws.Range("AR31").Value = 17
$ws.Range("AS31").Value = 15
$ws.Range("AT31").Value = 5
$ws.Range("AU31").Value = 17
$ws.Range("AV31").Value = 1
$ws.Range("AW31").Value = 8
$ws.Range("AX31").Value = 1
$ws.Range("AY31").Value = 1
$ws.Range("AZ31").Value = 6
$ws.Range("BA31").Value = 1
$ws.Range("BB31").Value = 9
$ws.Range("BC31").Value = 2
$ws.Range("BD31").Value = 7
$ws.Range("BE31").Value = 14
$ws.Range("BG31").Value = 27
$ws.Range("BH31").Value = 3
$ws.Range("BI31").Value = 54

# Row 32 (Apr 5)
$ws.Range("A32").Value = 43926
$ws.Range("C32").Value = 1326
$ws.Range("D32").Value = 10
$ws.Range("F32").Value = 15
$ws.Range("H32").Value = 226
$ws.Range("J32").Value = 286
$ws.Range("L32").Value = 237
$ws.Range("N32").Value = 222
$ws.Range("P32").Value = 162
$ws.Range("R32").Value = 99
$ws.Range("T32").Value = 66
$ws.Range("V32").Value = 3
$ws.Range("X32").Value = 609
$ws.Range("Y32").Value = 710
$ws.Range("Z32").Value = 7
$ws.Range("AA32").Value = 249
$ws.Range("AB32").Value = 94
$ws.Range("AC32").Value = 19
$ws.Range("AD32").Value = 41
$ws.Range("AE32").Value = 95
$ws.Range("AF32").Value = 4
$ws.Range("AG32").Value = 8
$ws.Range("AH32").Value = 72
$ws.Range("AI32").Value = 29
$ws.Range("AJ32").Value = 28
$ws.Range("AK32").Value = 6
$ws.Range("AL32").Value = 18
$ws.Range("AM32").Value = 12
$ws.Range("AN32").Value = 19
$ws.Range("AO32").Value = 30
$ws.Range("AP32").Value = 13
$ws.Range("AQ32").Value = 686
$ws.Range("AR32").Value = 17
$ws.Range("AS32").Value = 16
$ws.Range("AT32").Value = 5
$ws.Range("AU32").Value = 18
$ws.Range("AV32").Value = 1
$ws.Range("AW32").Value = 8
$ws.Range("AX32").Value = 1
$ws.Range("AY32").Value = 1
$ws.Range("AZ32").Value = 6
$ws.Range("BA32").Value = 1
$ws.Range("BB32").Value = 10
$ws.Range("BC32").Value = 2
$ws.Range("BD32").Value = 7
$ws.Range("BE32").Value = 14
$ws.Range("BG32").Value = 28
$ws.Range("BH32").Value = 3
$ws.Range("BI32").Value = 57

# Update view state to match
$ws.Range("BJ32").Select()
